$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 29-30; this pushes the existing rows 29-45
# down to 31-47 and keeps their contents unchanged.
$ws.Rows("29:30").Insert()

# New row 29 - weekly "Primera" quality entry
$ws.Cells.Item(29, 1).Value = 1
$ws.Cells.Item(29, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(29, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(29, 4).Value = 44767
$ws.Cells.Item(29, 5).Value = 15
$ws.Cells.Item(29, 6).Value = "Fruta"
$ws.Cells.Item(29, 7).Value = 100108
$ws.Cells.Item(29, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(29, 9).Value = 100108001
$ws.Cells.Item(29, 10).Value = "Guayaba"
$ws.Cells.Item(29, 11).Value = "Sin especificar"
$ws.Cells.Item(29, 12).Value = "Primera"
$ws.Cells.Item(29, 13).Value = 200
$ws.Cells.Item(29, 14).Value = 800
$ws.Cells.Item(29, 15).Value = 900
$ws.Cells.Item(29, 16).Value = 850
$ws.Cells.Item(29, 17).Value = "$/kilo (en caja de 10 kilos )"
$ws.Cells.Item(29, 18).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(29, 19).Value = 850
$ws.Cells.Item(29, 20).Value = 1

# New row 30 - weekly "Segunda" quality entry
$ws.Cells.Item(30, 1).Value = 1
$ws.Cells.Item(30, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(30, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(30, 4).Value = 44767
$ws.Cells.Item(30, 5).Value = 15
$ws.Cells.Item(30, 6).Value = "Fruta"
$ws.Cells.Item(30, 7).Value = 100108
$ws.Cells.Item(30, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(30, 9).Value = 100108001
$ws.Cells.Item(30, 10).Value = "Guayaba"
$ws.Cells.Item(30, 11).Value = "Sin especificar"
$ws.Cells.Item(30, 12).Value = "Segunda"
$ws.Cells.Item(30, 13).Value = 200
$ws.Cells.Item(30, 14).Value = 600
$ws.Cells.Item(30, 15).Value = 700
$ws.Cells.Item(30, 16).Value = 650
$ws.Cells.Item(30, 17).Value = "$/kilo (en caja de 10 kilos )"
$ws.Cells.Item(30, 18).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(30, 19).Value = 650
$ws.Cells.Item(30, 20).Value = 1
